$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new centered/underlined title paragraph and a following
#    blank paragraph at the very top of the document, ahead of the
#    existing "Overview of Project" paragraph.
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$insertionRange = $firstPara.Range
$insertionRange.Collapse(1)

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="NormalWeb"/>' + `
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' + `
    '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="240" w:afterAutospacing="0"/>' + `
    '<w:jc w:val="center"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>' + `
      '<w:b/>' + `
      '<w:bCs/>' + `
      '<w:color w:val="24292F"/>' + `
      '<w:sz w:val="28"/>' + `
      '<w:szCs w:val="28"/>' + `
      '<w:u w:val="single"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>' + `
      '<w:b/>' + `
      '<w:bCs/>' + `
      '<w:color w:val="24292F"/>' + `
      '<w:sz w:val="28"/>' + `
      '<w:szCs w:val="28"/>' + `
      '<w:u w:val="single"/>' + `
    '</w:rPr>' + `
    '<w:t>Kickstarter Analysis Written Assignment</w:t>' + `
  '</w:r>' + `
'</w:p>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="NormalWeb"/>' + `
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' + `
    '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="240" w:afterAutospacing="0"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>' + `
      '<w:b/>' + `
      '<w:bCs/>' + `
      '<w:color w:val="24292F"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
'</w:p>' + `
'</w:body>' + `
'</w:document>' + `
'</pkg:xmlData>' + `
'</pkg:part>' + `
'</pkg:package>'

$insertionRange.InsertXML($titleXml)

# ------------------------------------------------------------------
# 2) Drop the stray <w:lastRenderedPageBreak/> hint that sits in front
#    of the "Figure 2" caption run -- re-typing the run through Find /
#    Replace regenerates the run without the rendering-only element.
#    Locate the lone "Figure 2" caption paragraph (centered, bold,
#    exact text "Figure 2") rather than any of the inline mentions of
#    "Figure 2" inside body text elsewhere in the document.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Figure 2`r") {
        $capRange = $p.Range
        $capRange.Find.Execute("Figure 2", $true, $true, $false, $false, $false, `
                                $true, 1, $false, "Figure 2", 2)
        break
    }
}
